# Update data dictionary export/import: use "display" values for data types
# (Date / Plain / Multiple Choice) instead of the raw internal keys, on both
# the "caseType1" and "caseType2" sheets (C2:C4 = Data Type column).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("caseType1")
$ws2 = $wb.Worksheets.Item("caseType1-vl")
$ws3 = $wb.Worksheets.Item("caseType2")
$ws4 = $wb.Worksheets.Item("caseType2-vl")

$ws1.Range("C2").Value = "Date"
$ws1.Range("C3").Value = "Plain"
$ws1.Range("C4").Value = "Multiple Choice"

$ws3.Range("C2").Value = "Date"
$ws3.Range("C3").Value = "Plain"
$ws3.Range("C4").Value = "Multiple Choice"

# Column widths were re-fitted (e.g. after an "optimal width" pass) once the
# Data Type text changed. ColumnWidth is expressed in characters; set the
# explicit widths for the affected columns on each sheet.
$ws1.Columns.Item(1).ColumnWidth = 19.333333333333336
$ws1.Columns.Item(2).ColumnWidth = 24.166666666666668
$ws1.Columns.Item(3).ColumnWidth = 25.333333333333336

$ws2.Columns.Item(1).ColumnWidth = 17.833333333333336
$ws2.Columns.Item(2).ColumnWidth = 30.833333333333336
$ws2.Columns.Item(3).ColumnWidth = 28.5

$ws4.Columns.Item(1).ColumnWidth = 35.666666666666664
$ws4.Columns.Item(2).ColumnWidth = 31.666666666666668
$ws4.Columns.Item(3).ColumnWidth = 37.0

# Move the cached selection on "caseType2" to D4, then restore the
# originally-active sheet ("caseType2-vl") so the active tab is unchanged.
$ws3.Range("D4").Select()
$ws4.Activate()
